$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Target cluster was "MuSCs" (old rows 8, 9, 10)
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Update remaining data rows (2-7) with refreshed TPM-derived values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.324075666666667
$ws.Range("H2").Value = 3.972227
$ws.Range("I2").Value = 0.01675578032580584
$ws.Range("J2").Value = 0.01684165790066494
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1074926666666667
$ws.Range("N2").Value = 0.322478
$ws.Range("O2").Value = 0.01930181557781338
$ws.Range("P2").Value = 0.01930181557781338
$ws.Range("Q2").Value = 0.1423284242784444
$ws.Range("R2").Value = 1.280955818506
$ws.Range("S2").Value = 0.0003234169817110581
$ws.Range("T2").Value = 0.0003250745748232583

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.324075666666667
$ws.Range("H3").Value = 3.972227
$ws.Range("I3").Value = 0.01675578032580584
$ws.Range("J3").Value = 0.01684165790066494
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.461551666666666
$ws.Range("N3").Value = 16.384655
$ws.Range("O3").Value = 0.9806981844221867
$ws.Range("P3").Value = 0.9806981844221866
$ws.Range("Q3").Value = 7.23150766407611
$ws.Range("R3").Value = 65.083568976685
$ws.Range("S3").Value = 0.01643236334409478
$ws.Range("T3").Value = 0.01651658332584168

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 76.48912033333333
$ws.Range("H4").Value = 229.467361
$ws.Range("I4").Value = 0.9679468703219594
$ws.Range("J4").Value = 0.9729078406975189
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1074926666666667
$ws.Range("N4").Value = 0.322478
$ws.Range("O4").Value = 0.01930181557781338
$ws.Range("P4").Value = 0.01930181557781338
$ws.Range("Q4").Value = 8.222019515617555
$ws.Range("R4").Value = 73.99817564055799
$ws.Range("S4").Value = 0.01868313198007611
$ws.Range("T4").Value = 0.01877888771535215

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 76.48912033333333
$ws.Range("H5").Value = 229.467361
$ws.Range("I5").Value = 0.9679468703219594
$ws.Range("J5").Value = 0.9729078406975189
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.461551666666666
$ws.Range("N5").Value = 16.384655
$ws.Range("O5").Value = 0.9806981844221867
$ws.Range("P5").Value = 0.9806981844221866
$ws.Range("Q5").Value = 417.7492826383838
$ws.Range("R5").Value = 3759.743543745454
$ws.Range("S5").Value = 0.9492637383418834
$ws.Range("T5").Value = 0.9541289529821667

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.2088275
$ws.Range("H6").Value = 2.417655
$ws.Range("I6").Value = 0.0152973493522347
$ws.Range("J6").Value = 0.01025050140181618
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1074926666666667
$ws.Range("N6").Value = 0.322478
$ws.Range("O6").Value = 0.01930181557781338
$ws.Range("P6").Value = 0.01930181557781338
$ws.Range("Q6").Value = 0.129940091515
$ws.Range("R6").Value = 0.7796405490899999
$ws.Range("S6").Value = 0.0002952666160262172
$ws.Range("T6").Value = 0.0001978532876379735

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.2088275
$ws.Range("H7").Value = 2.417655
$ws.Range("I7").Value = 0.0152973493522347
$ws.Range("J7").Value = 0.01025050140181618
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.461551666666666
$ws.Range("N7").Value = 16.384655
$ws.Range("O7").Value = 0.9806981844221867
$ws.Range("P7").Value = 0.9806981844221866
$ws.Range("Q7").Value = 6.602073847337499
$ws.Range("R7").Value = 39.612443084025
$ws.Range("S7").Value = 0.01500208273620848
$ws.Range("T7").Value = 0.01005264811417821
